$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("front-end")

# --- Value changes on the front-end sheet ---
# E9: "PM" -> "F"
$ws.Range("E9").Value = "F"
# E12: "U" -> "F"
$ws.Range("E12").Value = "F"
# G12: "Jossias" -> cleared
$ws.Range("G12").Value = ""
# E23: "EP" -> cleared
$ws.Range("E23").Value = ""
# G23: "JOSSIAS" -> cleared
$ws.Range("G23").Value = ""

# --- View / selection changes ---
# Move the active selection from G16 to F14 (also clears the stale
# topLeftCell="A8" scroll-anchor left over from before).
$ws.Range("F14").Select() | Out-Null
